$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates (testing results for remaining rows) ---
$ws.Range("C6").Value = "yes"
$ws.Range("G6").ClearContents()

$ws.Range("A50").Value = "yes"
$ws.Range("C50").Value = "yes"

$ws.Range("A51").Value = "yes"
$ws.Range("C51").Value = "n/a"

$ws.Range("A52").Value = "yes"
$ws.Range("C52").Value = "yes"
$ws.Range("F52").Value = "BH recommends users be asked if they wish to open Python code."

$ws.Range("A53").Value = "yes"
$ws.Range("C53").Value = "yes"

# Updated note text for the "will not convert" row
$ws.Range("F56").Value = "Cannot not convert zip to TRZ with mp4  video"

# --- New column F formatting (Java Notes column) ---
$ws.Columns("F").ColumnWidth = 15.5

# --- Add hyperlinks to the four problem-report entries ---
$ws.Hyperlinks.Add($ws.Range("B8"), "https://www.compadre.org/osp/document/ServeFile.cfm?ID=14630", "", "", $ws.Range("B8").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B17"), "https://www.compadre.org/osp/document/ServeFile.cfm?ID=14630", "", "", $ws.Range("B17").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B28"), "https://www.compadre.org/osp/document/ServeFile.cfm?ID=14630", "", "", $ws.Range("B28").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B56"), "https://www.compadre.org/osp/document/ServeFile.cfm?ID=14630", "", "", $ws.Range("B56").Value2) | Out-Null

# --- Restore view state: freeze header row, scroll to row 7, select B8 ---
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B8").Select()
